$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$__style = $ws.Cells.Item(2, 4).Style
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '67.731.33'
$ws.Cells.Item(2, 4).Style = $__style
$ws.Cells.Item(2, 5).Value = '  -0.78%  '

$__style = $ws.Cells.Item(3, 4).Style
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '3.534.74'
$ws.Cells.Item(3, 4).Style = $__style
$ws.Cells.Item(3, 5).Value = '  -2.81%  '

$ws.Cells.Item(4, 5).Value = '  -0.02%  '

$__style = $ws.Cells.Item(5, 4).Style
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '203.31'
$ws.Cells.Item(5, 4).Style = $__style
$ws.Cells.Item(5, 5).Value = '  +2.97%  '

$__style = $ws.Cells.Item(6, 4).Style
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '555.99'
$ws.Cells.Item(6, 4).Style = $__style
$ws.Cells.Item(6, 5).Value = '  -3.91%  '

$__style = $ws.Cells.Item(7, 4).Style
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '3.525.65'
$ws.Cells.Item(7, 4).Style = $__style
$ws.Cells.Item(7, 5).Value = '  -2.93%  '

$ws.Cells.Item(8, 5).Value = '  -1.72%  '

$ws.Cells.Item(9, 5).Value = '  -0.04%  '

$ws.Cells.Item(10, 2).Value = 'Cardano'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$__style = $ws.Cells.Item(10, 4).Style
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.659'
$ws.Cells.Item(10, 4).Style = $__style
$ws.Cells.Item(10, 5).Value = '  -2.67%  '

$ws.Cells.Item(11, 2).Value = 'Avalanche'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$__style = $ws.Cells.Item(11, 4).Style
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '63.95'
$ws.Cells.Item(11, 4).Style = $__style
$ws.Cells.Item(11, 5).Value = '  +12.70%  '

$ws.Cells.Item(12, 5).Value = '  -7.07%  '

$__style = $ws.Cells.Item(13, 4).Style
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.0000272'
$ws.Cells.Item(13, 4).Style = $__style
$ws.Cells.Item(13, 5).Value = '  -7.38%  '

$__style = $ws.Cells.Item(14, 4).Style
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '9.94'
$ws.Cells.Item(14, 4).Style = $__style
$ws.Cells.Item(14, 5).Value = '  -1.46%  '

$__style = $ws.Cells.Item(15, 4).Style
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '4.097.51'
$ws.Cells.Item(15, 4).Style = $__style
$ws.Cells.Item(15, 5).Value = '  -2.97%  '

$__style = $ws.Cells.Item(16, 4).Style
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '3.536.28'
$ws.Cells.Item(16, 4).Style = $__style
$ws.Cells.Item(16, 5).Value = '  -2.96%  '

$ws.Cells.Item(17, 5).Value = '  -1.89%  '

$__style = $ws.Cells.Item(18, 4).Style
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '18.53'
$ws.Cells.Item(18, 4).Style = $__style
$ws.Cells.Item(18, 5).Value = '  -0.48%  '

$__style = $ws.Cells.Item(19, 4).Style
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '67.533.21'
$ws.Cells.Item(19, 4).Style = $__style
$ws.Cells.Item(19, 5).Value = '  -1.07%  '

$__style = $ws.Cells.Item(20, 4).Style
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '11.92'
$ws.Cells.Item(20, 4).Style = $__style
$ws.Cells.Item(20, 5).Value = '  -5.22%  '

$ws.Cells.Item(21, 5).Value = '  -5.18%  '

$__style = $ws.Cells.Item(22, 4).Style
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '394.88'
$ws.Cells.Item(22, 4).Style = $__style
$ws.Cells.Item(22, 5).Value = '  -2.12%  '

$__style = $ws.Cells.Item(23, 4).Style
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '12.11'
$ws.Cells.Item(23, 4).Style = $__style
$ws.Cells.Item(23, 5).Value = '  -8.23%  '

$ws.Cells.Item(24, 5).Value = '  -5.65%  '

$__style = $ws.Cells.Item(25, 4).Style
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '83.44'
$ws.Cells.Item(25, 4).Style = $__style
$ws.Cells.Item(25, 5).Value = '  -3.08%  '

$__style = $ws.Cells.Item(26, 4).Style
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '3.93'
$ws.Cells.Item(26, 4).Style = $__style
$ws.Cells.Item(26, 5).Value = '  +1.35%  '

$__style = $ws.Cells.Item(27, 4).Style
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '2.85'
$ws.Cells.Item(27, 4).Style = $__style
$ws.Cells.Item(27, 5).Value = '  -4.03%  '

$__style = $ws.Cells.Item(28, 4).Style
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '12.33'
$ws.Cells.Item(28, 4).Style = $__style
$ws.Cells.Item(28, 5).Value = '  -2.67%  '

$__style = $ws.Cells.Item(29, 4).Style
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '8.93'
$ws.Cells.Item(29, 4).Style = $__style
$ws.Cells.Item(29, 5).Value = '  -3.10%  '

$__style = $ws.Cells.Item(30, 4).Style
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '727.30'
$ws.Cells.Item(30, 4).Style = $__style
$ws.Cells.Item(30, 5).Value = '  +5.52%  '

$ws.Cells.Item(31, 5).Value = '  -1.95%  '

$__style = $ws.Cells.Item(32, 4).Style
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '7.18'
$ws.Cells.Item(32, 4).Style = $__style
$ws.Cells.Item(32, 5).Value = '  -13.11%  '

$__style = $ws.Cells.Item(33, 4).Style
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '11.79'
$ws.Cells.Item(33, 4).Style = $__style
$ws.Cells.Item(33, 5).Value = '  -3.75%  '

$__style = $ws.Cells.Item(34, 4).Style
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '64.30'
$ws.Cells.Item(34, 4).Style = $__style
$ws.Cells.Item(34, 5).Value = '  -0.60%  '

$__style = $ws.Cells.Item(35, 4).Style
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.112'
$ws.Cells.Item(35, 4).Style = $__style
$ws.Cells.Item(35, 5).Value = '  -4.70%  '

$__style = $ws.Cells.Item(36, 4).Style
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '38.94'
$ws.Cells.Item(36, 4).Style = $__style
$ws.Cells.Item(36, 5).Value = '  -9.06%  '

$ws.Cells.Item(37, 5).Value = '  -6.48%  '

$ws.Cells.Item(38, 5).Value = '  +0.08%  '

$ws.Cells.Item(39, 2).Value = 'ThetaToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$__style = $ws.Cells.Item(39, 4).Style
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '3.10'
$ws.Cells.Item(39, 4).Style = $__style
$ws.Cells.Item(39, 5).Value = '  -2.31%  '

$ws.Cells.Item(40, 2).Value = 'Kaspa'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$__style = $ws.Cells.Item(40, 4).Style
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.133'
$ws.Cells.Item(40, 4).Style = $__style
$ws.Cells.Item(40, 5).Value = '  -5.04%  '

$ws.Cells.Item(41, 2).Value = 'Maker'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$__style = $ws.Cells.Item(41, 4).Style
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '3.080.91'
$ws.Cells.Item(41, 4).Style = $__style
$ws.Cells.Item(41, 5).Value = '  -4.40%  '

$ws.Cells.Item(42, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$__style = $ws.Cells.Item(42, 4).Style
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.998'
$ws.Cells.Item(42, 4).Style = $__style
$ws.Cells.Item(42, 5).Value = '  -0.08%  '

$__style = $ws.Cells.Item(43, 4).Style
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.0₃0688'
$ws.Cells.Item(43, 4).Style = $__style
$ws.Cells.Item(43, 5).Value = '  -13.03%  '

$ws.Cells.Item(44, 5).Value = '  -11.22%  '

$ws.Cells.Item(45, 2).Value = 'VeChain'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$__style = $ws.Cells.Item(45, 4).Style
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.0414'
$ws.Cells.Item(45, 4).Style = $__style
$ws.Cells.Item(45, 5).Value = '  -1.98%  '

$ws.Cells.Item(46, 2).Value = 'dogwifhat'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(46, 5).Value = '  -10.43%  '

$ws.Cells.Item(47, 2).Value = 'WEMIXToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$__style = $ws.Cells.Item(47, 4).Style
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '2.76'
$ws.Cells.Item(47, 4).Style = $__style
$ws.Cells.Item(47, 5).Value = '  +5.12%  '

$ws.Cells.Item(48, 5).Value = '  -2.93%  '

$__style = $ws.Cells.Item(49, 4).Style
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '138.43'
$ws.Cells.Item(49, 4).Style = $__style
$ws.Cells.Item(49, 5).Value = '  -3.50%  '

$__style = $ws.Cells.Item(50, 4).Style
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '8.31'
$ws.Cells.Item(50, 4).Style = $__style
$ws.Cells.Item(50, 5).Value = '  -7.27%  '

$__style = $ws.Cells.Item(51, 4).Style
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '2.88'
$ws.Cells.Item(51, 4).Style = $__style
$ws.Cells.Item(51, 5).Value = '  -7.37%  '
